# Updated jan-feb 2023 rates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 86 -> year 2023, month "Jan" (shared string index 18)
$ws.Range("C86").Value = 3.34
$ws.Range("D86").Value = 4.12
$ws.Range("E86").Value = 3.94
$ws.Range("F86").Value = 4.19
$ws.Range("G86").Value = 3.48

# Row 87 -> year 2023, month "Feb" (shared string index 17)
$ws.Range("C87").Value = 3.5
$ws.Range("D87").Value = 4.3
$ws.Range("E87").Value = 4.19
$ws.Range("F87").Value = 3.99
$ws.Range("G87").Value = 3.91

# Update the active selection on the frozen (bottom-left) pane to H82
$ws.Range("H82").Select()
